$wb = $excel.ActiveWorkbook

$sheetChanges = @{
    "ALC" = @(
        @{ Cell = "H2"; Value = 300 },
        @{ Cell = "J2"; Value = 0 },
        @{ Cell = "L2"; Value = 0 },
        @{ Cell = "N2"; Value = $null },
        @{ Cell = "H9"; Value = 1337254.9 },
        @{ Cell = "J9"; Value = 497.5 },
        @{ Cell = "L9"; Value = 497.5 },
        @{ Cell = "N9"; Value = -835.5 },
        @{ Cell = "H32"; Value = 20834908 },
        @{ Cell = "J32"; Value = 2900 },
        @{ Cell = "L32"; Value = 2900 },
        @{ Cell = "N32"; Value = -3552 },
        @{ Cell = "H34"; Value = 1316.375 },
        @{ Cell = "I34"; Value = 1316.375 },
        @{ Cell = "J34"; Value = 0 },
        @{ Cell = "K34"; Value = 1316.375 },
        @{ Cell = "L34"; Value = 0 },
        @{ Cell = "M34"; Value = -1113.375 },
        @{ Cell = "N34"; Value = $null },
        @{ Cell = "H36"; Value = 1316.375 },
        @{ Cell = "I36"; Value = 1316.375 },
        @{ Cell = "J36"; Value = 0 },
        @{ Cell = "K36"; Value = 1316.375 },
        @{ Cell = "L36"; Value = 0 },
        @{ Cell = "M36"; Value = -601.375 },
        @{ Cell = "N36"; Value = $null },
        @{ Cell = "H43"; Value = 1902.5714 },
        @{ Cell = "I43"; Value = 1279 },
        @{ Cell = "J43"; Value = 2006.5 },
        @{ Cell = "K43"; Value = 1279 },
        @{ Cell = "L43"; Value = 2006.5 },
        @{ Cell = "M43"; Value = -1210 },
        @{ Cell = "N43"; Value = -2144.5 },
        @{ Cell = "H61"; Value = 3671 },
        @{ Cell = "I61"; Value = 3671 },
        @{ Cell = "K61"; Value = 11013 },
        @{ Cell = "M61"; Value = -10841 },
        @{ Cell = "H74"; Value = 12739.35 },
        @{ Cell = "I74"; Value = 12739.35 },
        @{ Cell = "K74"; Value = 12739.35 },
        @{ Cell = "M74"; Value = -11803.35 },
        @{ Cell = "H76"; Value = 3714 },
        @{ Cell = "I76"; Value = 3714 },
        @{ Cell = "J76"; Value = 0 },
        @{ Cell = "K76"; Value = 3714 },
        @{ Cell = "L76"; Value = 0 },
        @{ Cell = "M76"; Value = -3399 },
        @{ Cell = "N76"; Value = $null },
        @{ Cell = "H77"; Value = 12739.35 },
        @{ Cell = "I77"; Value = 12739.35 },
        @{ Cell = "K77"; Value = 63696.75 },
        @{ Cell = "M77"; Value = -59016.75 },
        @{ Cell = "H79"; Value = 3714 },
        @{ Cell = "I79"; Value = 3714 },
        @{ Cell = "J79"; Value = 0 },
        @{ Cell = "K79"; Value = 3714 },
        @{ Cell = "L79"; Value = 0 },
        @{ Cell = "M79"; Value = -2622 },
        @{ Cell = "N79"; Value = $null },
        @{ Cell = "H86"; Value = 142860500 },
        @{ Cell = "I86"; Value = 1000000000 },
        @{ Cell = "J86"; Value = 3916.6667 },
        @{ Cell = "K86"; Value = 1000000000 },
        @{ Cell = "L86"; Value = 3916.6667 },
        @{ Cell = "M86"; Value = -999998877 },
        @{ Cell = "N86"; Value = -6162.6667 },
        @{ Cell = "H88"; Value = 9095318 },
        @{ Cell = "I88"; Value = 33336000 },
        @{ Cell = "J88"; Value = 5062 },
        @{ Cell = "K88"; Value = 33336000 },
        @{ Cell = "L88"; Value = 5062 },
        @{ Cell = "M88"; Value = -33335594 },
        @{ Cell = "N88"; Value = -5874 },
        @{ Cell = "H89"; Value = 142860500 },
        @{ Cell = "I89"; Value = 1000000000 },
        @{ Cell = "J89"; Value = 3916.6667 },
        @{ Cell = "K89"; Value = 5000000000 },
        @{ Cell = "L89"; Value = 19583.3335 },
        @{ Cell = "M89"; Value = -4999994384 },
        @{ Cell = "N89"; Value = -30815.3335 },
        @{ Cell = "H91"; Value = 9095318 },
        @{ Cell = "I91"; Value = 33336000 },
        @{ Cell = "J91"; Value = 5062 },
        @{ Cell = "K91"; Value = 33336000 },
        @{ Cell = "L91"; Value = 5062 },
        @{ Cell = "M91"; Value = -33334596 },
        @{ Cell = "N91"; Value = -7870 },
        @{ Cell = "H92"; Value = 1974.7778 },
        @{ Cell = "J92"; Value = 1649.75 },
        @{ Cell = "L92"; Value = 1649.75 },
        @{ Cell = "N92"; Value = -4145.75 },
        @{ Cell = "H99"; Value = 311.83334 },
        @{ Cell = "I99"; Value = 311.83334 },
        @{ Cell = "K99"; Value = 935.5000200000001 },
        @{ Cell = "M99"; Value = 562.4999799999999 },
        @{ Cell = "H103"; Value = 5103511.5 },
        @{ Cell = "I103"; Value = 2026.5555 },
        @{ Cell = "J103"; Value = 14286184 },
        @{ Cell = "K103"; Value = 6079.666499999999 },
        @{ Cell = "L103"; Value = 42858552 },
        @{ Cell = "M103"; Value = -5493.666499999999 },
        @{ Cell = "N103"; Value = -42859724 },
        @{ Cell = "H107"; Value = 364.875 },
        @{ Cell = "I107"; Value = 364.875 },
        @{ Cell = "K107"; Value = 364.875 },
        @{ Cell = "M107"; Value = 1555.125 },
        @{ Cell = "H113"; Value = 4918.4443 },
        @{ Cell = "I113"; Value = 4459 },
        @{ Cell = "J113"; Value = 5837.3335 },
        @{ Cell = "K113"; Value = 4459 },
        @{ Cell = "L113"; Value = 5837.3335 },
        @{ Cell = "M113"; Value = -1205 },
        @{ Cell = "N113"; Value = -12345.3335 },
        @{ Cell = "H117"; Value = 0 },
        @{ Cell = "J117"; Value = 0 },
        @{ Cell = "L117"; Value = 0 },
        @{ Cell = "N117"; Value = $null },
        @{ Cell = "H125"; Value = 999.875 },
        @{ Cell = "I125"; Value = 999.6667 },
        @{ Cell = "K125"; Value = 8997.0003 },
        @{ Cell = "M125"; Value = -6537.0003 },
        @{ Cell = "H127"; Value = 474.8 },
        @{ Cell = "I127"; Value = 307.375 },
        @{ Cell = "K127"; Value = 922.125 },
        @{ Cell = "M127"; Value = 4037.875 },
        @{ Cell = "H128"; Value = 0 },
        @{ Cell = "J128"; Value = 0 },
        @{ Cell = "L128"; Value = 0 },
        @{ Cell = "N128"; Value = $null },
        @{ Cell = "H129"; Value = 866.5 },
        @{ Cell = "I129"; Value = 499 },
        @{ Cell = "J129"; Value = 1969 },
        @{ Cell = "K129"; Value = 1497 },
        @{ Cell = "L129"; Value = 5907 },
        @{ Cell = "M129"; Value = 3503 },
        @{ Cell = "N129"; Value = -15907 },
        @{ Cell = "H132"; Value = 8239.611000000001 },
        @{ Cell = "I132"; Value = 8665.471 },
        @{ Cell = "K132"; Value = 25996.413 },
        @{ Cell = "M132"; Value = -23466.413 },
        @{ Cell = "H135"; Value = 645.75 },
        @{ Cell = "I135"; Value = 329.70587 },
        @{ Cell = "J135"; Value = 2436.6667 },
        @{ Cell = "K135"; Value = 2967.35283 },
        @{ Cell = "L135"; Value = 21930.0003 },
        @{ Cell = "M135"; Value = -432.3528299999998 },
        @{ Cell = "N135"; Value = -27000.0003 },
        @{ Cell = "H137"; Value = 1568311.8 },
        @{ Cell = "I137"; Value = 2175705 },
        @{ Cell = "K137"; Value = 6527115 },
        @{ Cell = "M137"; Value = -6524565 },
        @{ Cell = "H138"; Value = 2665.2534 },
        @{ Cell = "I138"; Value = 1255 },
        @{ Cell = "J138"; Value = 2819.5 },
        @{ Cell = "K138"; Value = 3765 },
        @{ Cell = "L138"; Value = 8458.5 },
        @{ Cell = "M138"; Value = 1375 },
        @{ Cell = "N138"; Value = -18738.5 },
        @{ Cell = "H141"; Value = 3342 },
        @{ Cell = "I141"; Value = 3342 },
        @{ Cell = "K141"; Value = 10026 },
        @{ Cell = "M141"; Value = -4846 }
    )
    "ARM" = @(
        @{ Cell = "H6"; Value = 3897 },
        @{ Cell = "I6"; Value = 346 },
        @{ Cell = "K6"; Value = 346 },
        @{ Cell = "M6"; Value = -173 },
        @{ Cell = "H19"; Value = 4004 },
        @{ Cell = "I19"; Value = 3008 },
        @{ Cell = "J19"; Value = 5000 },
        @{ Cell = "K19"; Value = 3008 },
        @{ Cell = "L19"; Value = 5000 },
        @{ Cell = "M19"; Value = -2779 },
        @{ Cell = "N19"; Value = -5458 },
        @{ Cell = "H32"; Value = 2434.976 },
        @{ Cell = "I32"; Value = 2421.244 },
        @{ Cell = "K32"; Value = 2421.244 },
        @{ Cell = "M32"; Value = -2134.244 },
        @{ Cell = "H45"; Value = 34978.92 },
        @{ Cell = "I45"; Value = 72288.664 },
        @{ Cell = "J45"; Value = 2999.1428 },
        @{ Cell = "K45"; Value = 72288.664 },
        @{ Cell = "L45"; Value = 2999.1428 },
        @{ Cell = "M45"; Value = -71911.664 },
        @{ Cell = "N45"; Value = -3753.1428 },
        @{ Cell = "H61"; Value = 3886.3 },
        @{ Cell = "I61"; Value = 2099 },
        @{ Cell = "J61"; Value = 4084.889 },
        @{ Cell = "K61"; Value = 2099 },
        @{ Cell = "L61"; Value = 4084.889 },
        @{ Cell = "M61"; Value = -1887 },
        @{ Cell = "N61"; Value = -4508.889 },
        @{ Cell = "H63"; Value = 1000 },
        @{ Cell = "J63"; Value = 0 },
        @{ Cell = "L63"; Value = 0 },
        @{ Cell = "N63"; Value = $null },
        @{ Cell = "H66"; Value = 1000 },
        @{ Cell = "J66"; Value = 0 },
        @{ Cell = "L66"; Value = 0 },
        @{ Cell = "N66"; Value = $null },
        @{ Cell = "H74"; Value = 224385.8 },
        @{ Cell = "I74"; Value = 243645.44 },
        @{ Cell = "K74"; Value = 243645.44 },
        @{ Cell = "M74"; Value = -242771.44 },
        @{ Cell = "H77"; Value = 224385.8 },
        @{ Cell = "I77"; Value = 243645.44 },
        @{ Cell = "K77"; Value = 1218227.2 },
        @{ Cell = "M77"; Value = -1213859.2 },
        @{ Cell = "H86"; Value = 70000 },
        @{ Cell = "J86"; Value = 70000 },
        @{ Cell = "L86"; Value = 70000 },
        @{ Cell = "N86"; Value = -72372 },
        @{ Cell = "H89"; Value = 70000 },
        @{ Cell = "J89"; Value = 70000 },
        @{ Cell = "L89"; Value = 210000 },
        @{ Cell = "N89"; Value = -221856 },
        @{ Cell = "H97"; Value = 1299 },
        @{ Cell = "I97"; Value = 1299 },
        @{ Cell = "J97"; Value = 0 },
        @{ Cell = "K97"; Value = 1299 },
        @{ Cell = "L97"; Value = 0 },
        @{ Cell = "M97"; Value = -803 },
        @{ Cell = "N97"; Value = $null },
        @{ Cell = "H101"; Value = 0 },
        @{ Cell = "J101"; Value = 0 },
        @{ Cell = "L101"; Value = 0 },
        @{ Cell = "N101"; Value = $null },
        @{ Cell = "H102"; Value = 2559.4 },
        @{ Cell = "I102"; Value = 2145.4614 },
        @{ Cell = "K102"; Value = 2145.4614 },
        @{ Cell = "M102"; Value = -523.4614000000001 },
        @{ Cell = "H110"; Value = 9191.375 },
        @{ Cell = "I110"; Value = 9074 },
        @{ Cell = "K110"; Value = 9074 },
        @{ Cell = "M110"; Value = -7029 },
        @{ Cell = "H132"; Value = 1743.4667 },
        @{ Cell = "I132"; Value = 1096.7 },
        @{ Cell = "K132"; Value = 3290.1 },
        @{ Cell = "M132"; Value = -760.1000000000004 },
        @{ Cell = "H135"; Value = 88386 },
        @{ Cell = "J135"; Value = 108995 },
        @{ Cell = "L135"; Value = 108995 },
        @{ Cell = "N135"; Value = -119135 },
        @{ Cell = "H136"; Value = 3886.3 },
        @{ Cell = "I136"; Value = 2099 },
        @{ Cell = "J136"; Value = 4084.889 },
        @{ Cell = "K136"; Value = 6297 },
        @{ Cell = "L136"; Value = 12254.667 },
        @{ Cell = "M136"; Value = -3747 },
        @{ Cell = "N136"; Value = -17354.667 }
    )
    "BSM" = @(
        @{ Cell = "H6"; Value = 68999 },
        @{ Cell = "J6"; Value = 68999 },
        @{ Cell = "L6"; Value = 68999 },
        @{ Cell = "N6"; Value = -69225 },
        @{ Cell = "H20"; Value = 49024412 },
        @{ Cell = "I20"; Value = 52088188 },
        @{ Cell = "J20"; Value = 4009 },
        @{ Cell = "K20"; Value = 52088188 },
        @{ Cell = "L20"; Value = 4009 },
        @{ Cell = "M20"; Value = -52087941 },
        @{ Cell = "N20"; Value = -4503 },
        @{ Cell = "H86"; Value = 2354.348 },
        @{ Cell = "I86"; Value = 2226.8096 },
        @{ Cell = "J86"; Value = 3693.5 },
        @{ Cell = "K86"; Value = 2226.8096 },
        @{ Cell = "L86"; Value = 3693.5 },
        @{ Cell = "M86"; Value = -1103.8096 },
        @{ Cell = "N86"; Value = -5939.5 },
        @{ Cell = "H88"; Value = 0 },
        @{ Cell = "J88"; Value = 0 },
        @{ Cell = "L88"; Value = 0 },
        @{ Cell = "N88"; Value = $null },
        @{ Cell = "H89"; Value = 2354.348 },
        @{ Cell = "I89"; Value = 2226.8096 },
        @{ Cell = "J89"; Value = 3693.5 },
        @{ Cell = "K89"; Value = 11134.048 },
        @{ Cell = "L89"; Value = 18467.5 },
        @{ Cell = "M89"; Value = -5518.048000000001 },
        @{ Cell = "N89"; Value = -29699.5 },
        @{ Cell = "H91"; Value = 0 },
        @{ Cell = "J91"; Value = 0 },
        @{ Cell = "L91"; Value = 0 },
        @{ Cell = "N91"; Value = $null },
        @{ Cell = "H94"; Value = 125007270 },
        @{ Cell = "I94"; Value = 200010640 },
        @{ Cell = "K94"; Value = 200010640 },
        @{ Cell = "M94"; Value = -200010189 },
        @{ Cell = "H134"; Value = 2351.6572 },
        @{ Cell = "I134"; Value = 1728.3334 },
        @{ Cell = "J134"; Value = 3711.6365 },
        @{ Cell = "K134"; Value = 5185.0002 },
        @{ Cell = "L134"; Value = 11134.9095 },
        @{ Cell = "M134"; Value = -2650.0002 },
        @{ Cell = "N134"; Value = -16204.9095 }
    )
    "CRP" = @(
        @{ Cell = "H7"; Value = 3112.9412 },
        @{ Cell = "I7"; Value = 3988 },
        @{ Cell = "J7"; Value = 269 },
        @{ Cell = "K7"; Value = 3988 },
        @{ Cell = "L7"; Value = 269 },
        @{ Cell = "M7"; Value = -3875 },
        @{ Cell = "N7"; Value = -495 },
        @{ Cell = "H19"; Value = 2224.077 },
        @{ Cell = "I19"; Value = 2008.2727 },
        @{ Cell = "J19"; Value = 3411 },
        @{ Cell = "K19"; Value = 2008.2727 },
        @{ Cell = "L19"; Value = 3411 },
        @{ Cell = "M19"; Value = -1838.2727 },
        @{ Cell = "N19"; Value = -3751 },
        @{ Cell = "H24"; Value = 2224.077 },
        @{ Cell = "I24"; Value = 2008.2727 },
        @{ Cell = "J24"; Value = 3411 },
        @{ Cell = "K24"; Value = 2008.2727 },
        @{ Cell = "L24"; Value = 3411 },
        @{ Cell = "M24"; Value = -1838.2727 },
        @{ Cell = "N24"; Value = -3751 },
        @{ Cell = "H31"; Value = 4497.927 },
        @{ Cell = "I31"; Value = 3249.2 },
        @{ Cell = "K31"; Value = 3249.2 },
        @{ Cell = "M31"; Value = -2954.2 },
        @{ Cell = "H34"; Value = 4497.927 },
        @{ Cell = "I34"; Value = 3249.2 },
        @{ Cell = "K34"; Value = 3249.2 },
        @{ Cell = "M34"; Value = -3047.2 },
        @{ Cell = "H58"; Value = 2615.4666 },
        @{ Cell = "I58"; Value = 1691.2222 },
        @{ Cell = "J58"; Value = 4001.8333 },
        @{ Cell = "K58"; Value = 1691.2222 },
        @{ Cell = "L58"; Value = 4001.8333 },
        @{ Cell = "M58"; Value = -1488.2222 },
        @{ Cell = "N58"; Value = -4407.8333 },
        @{ Cell = "H88"; Value = 0 },
        @{ Cell = "J88"; Value = 0 },
        @{ Cell = "L88"; Value = 0 },
        @{ Cell = "N88"; Value = $null },
        @{ Cell = "H91"; Value = 0 },
        @{ Cell = "J91"; Value = 0 },
        @{ Cell = "L91"; Value = 0 },
        @{ Cell = "N91"; Value = $null },
        @{ Cell = "H99"; Value = 5032.8335 },
        @{ Cell = "I99"; Value = 3100 },
        @{ Cell = "J99"; Value = 5999.25 },
        @{ Cell = "K99"; Value = 3100 },
        @{ Cell = "L99"; Value = 5999.25 },
        @{ Cell = "M99"; Value = -1602 },
        @{ Cell = "N99"; Value = -8995.25 },
        @{ Cell = "H105"; Value = 3237.375 },
        @{ Cell = "I105"; Value = 3722.25 },
        @{ Cell = "J105"; Value = 2752.5 },
        @{ Cell = "K105"; Value = 3722.25 },
        @{ Cell = "L105"; Value = 2752.5 },
        @{ Cell = "M105"; Value = -1975.25 },
        @{ Cell = "N105"; Value = -6246.5 },
        @{ Cell = "H106"; Value = 108788 },
        @{ Cell = "J106"; Value = 108788 },
        @{ Cell = "L106"; Value = 108788 },
        @{ Cell = "N106"; Value = -111312 },
        @{ Cell = "H121"; Value = 49999 },
        @{ Cell = "J121"; Value = 49999 },
        @{ Cell = "L121"; Value = 49999 },
        @{ Cell = "N121"; Value = -52619 },
        @{ Cell = "H126"; Value = 5032.8335 },
        @{ Cell = "I126"; Value = 3100 },
        @{ Cell = "J126"; Value = 5999.25 },
        @{ Cell = "K126"; Value = 9300 },
        @{ Cell = "L126"; Value = 17997.75 },
        @{ Cell = "M126"; Value = -6830 },
        @{ Cell = "N126"; Value = -22937.75 },
        @{ Cell = "H132"; Value = 13338481 },
        @{ Cell = "I132"; Value = 5048.8 },
        @{ Cell = "J132"; Value = 22227436 },
        @{ Cell = "K132"; Value = 15146.4 },
        @{ Cell = "L132"; Value = 66682308 },
        @{ Cell = "M132"; Value = -12616.4 },
        @{ Cell = "N132"; Value = -66687368 },
        @{ Cell = "H134"; Value = 4770.278 },
        @{ Cell = "I134"; Value = 5326.5 },
        @{ Cell = "J134"; Value = 2823.5 },
        @{ Cell = "K134"; Value = 15979.5 },
        @{ Cell = "L134"; Value = 8470.5 },
        @{ Cell = "M134"; Value = -13444.5 },
        @{ Cell = "N134"; Value = -13540.5 },
        @{ Cell = "H136"; Value = 2615.4666 },
        @{ Cell = "I136"; Value = 1691.2222 },
        @{ Cell = "J136"; Value = 4001.8333 },
        @{ Cell = "K136"; Value = 5073.6666 },
        @{ Cell = "L136"; Value = 12005.4999 },
        @{ Cell = "M136"; Value = -2523.6666 },
        @{ Cell = "N136"; Value = -17105.4999 },
        @{ Cell = "H137"; Value = 86398.2 },
        @{ Cell = "J137"; Value = 86398.2 },
        @{ Cell = "L137"; Value = 86398.2 },
        @{ Cell = "N137"; Value = -96598.2 }
    )
    "CUL" = @(
        @{ Cell = "H2"; Value = 852.3871 },
        @{ Cell = "I2"; Value = 106.8125 },
        @{ Cell = "J2"; Value = 1647.6666 },
        @{ Cell = "K2"; Value = 640.875 },
        @{ Cell = "L2"; Value = 9885.999599999999 },
        @{ Cell = "M2"; Value = -527.875 },
        @{ Cell = "N2"; Value = -10111.9996 },
        @{ Cell = "H12"; Value = 534.8570999999999 },
        @{ Cell = "J12"; Value = 545.6667 },
        @{ Cell = "L12"; Value = 1637.0001 },
        @{ Cell = "N12"; Value = -1983.0001 },
        @{ Cell = "H23"; Value = 6616.385 },
        @{ Cell = "J23"; Value = 7701.1816 },
        @{ Cell = "L23"; Value = 23103.5448 },
        @{ Cell = "N23"; Value = -23573.5448 },
        @{ Cell = "H38"; Value = 590.8333 },
        @{ Cell = "I38"; Value = 145.33333 },
        @{ Cell = "J38"; Value = 1036.3334 },
        @{ Cell = "K38"; Value = 435.99999 },
        @{ Cell = "L38"; Value = 3109.0002 },
        @{ Cell = "M38"; Value = -88.99998999999997 },
        @{ Cell = "N38"; Value = -3803.0002 },
        @{ Cell = "H44"; Value = 2869.125 },
        @{ Cell = "I44"; Value = 167.66667 },
        @{ Cell = "K44"; Value = 503.00001 },
        @{ Cell = "M44"; Value = -105.00001 },
        @{ Cell = "H55"; Value = 4826.7827 },
        @{ Cell = "J55"; Value = 5478.6665 },
        @{ Cell = "L55"; Value = 16435.9995 },
        @{ Cell = "N55"; Value = -16789.9995 },
        @{ Cell = "H68"; Value = 9100071 },
        @{ Cell = "J68"; Value = 16678049 },
        @{ Cell = "L68"; Value = 50034147 },
        @{ Cell = "N68"; Value = -50035769 },
        @{ Cell = "H71"; Value = 9100071 },
        @{ Cell = "J71"; Value = 16678049 },
        @{ Cell = "L71"; Value = 150102441 },
        @{ Cell = "N71"; Value = -150110553 },
        @{ Cell = "H81"; Value = 1303 },
        @{ Cell = "J81"; Value = 1462.25 },
        @{ Cell = "L81"; Value = 4386.75 },
        @{ Cell = "N81"; Value = -6632.75 },
        @{ Cell = "H84"; Value = 1303 },
        @{ Cell = "J84"; Value = 1462.25 },
        @{ Cell = "L84"; Value = 13160.25 },
        @{ Cell = "N84"; Value = -24392.25 },
        @{ Cell = "H88"; Value = 8833.333000000001 },
        @{ Cell = "J88"; Value = 8833.333000000001 },
        @{ Cell = "L88"; Value = 26499.999 },
        @{ Cell = "N88"; Value = -27355.999 },
        @{ Cell = "H91"; Value = 8833.333000000001 },
        @{ Cell = "J91"; Value = 8833.333000000001 },
        @{ Cell = "L91"; Value = 26499.999 },
        @{ Cell = "N91"; Value = -29463.999 },
        @{ Cell = "H92"; Value = 478.33334 },
        @{ Cell = "I92"; Value = 270 },
        @{ Cell = "K92"; Value = 810 },
        @{ Cell = "M92"; Value = 438 },
        @{ Cell = "H93"; Value = 0 },
        @{ Cell = "I93"; Value = 0 },
        @{ Cell = "K93"; Value = 0 },
        @{ Cell = "M93"; Value = $null },
        @{ Cell = "H97"; Value = 1003250 },
        @{ Cell = "I97"; Value = 2500175 },
        @{ Cell = "J97"; Value = 5300 },
        @{ Cell = "K97"; Value = 7500525 },
        @{ Cell = "L97"; Value = 15900 },
        @{ Cell = "M97"; Value = -7500029 },
        @{ Cell = "N97"; Value = -16892 },
        @{ Cell = "H98"; Value = 225 },
        @{ Cell = "J98"; Value = 100 },
        @{ Cell = "L98"; Value = 300 },
        @{ Cell = "N98"; Value = -3296 },
        @{ Cell = "H126"; Value = 0 },
        @{ Cell = "J126"; Value = 0 },
        @{ Cell = "L126"; Value = 0 },
        @{ Cell = "N126"; Value = $null },
        @{ Cell = "H129"; Value = 33334810 },
        @{ Cell = "I129"; Value = 0 },
        @{ Cell = "J129"; Value = 33334810 },
        @{ Cell = "K129"; Value = 0 },
        @{ Cell = "L129"; Value = 100004430 },
        @{ Cell = "M129"; Value = $null },
        @{ Cell = "N129"; Value = -100014430 },
        @{ Cell = "H131"; Value = 12502447 },
        @{ Cell = "I131"; Value = 62503452 },
        @{ Cell = "J131"; Value = 2196.5 },
        @{ Cell = "K131"; Value = 187510356 },
        @{ Cell = "L131"; Value = 6589.5 },
        @{ Cell = "M131"; Value = -187505316 },
        @{ Cell = "N131"; Value = -16669.5 },
        @{ Cell = "H139"; Value = 5531.2144 },
        @{ Cell = "I139"; Value = 3407.6 },
        @{ Cell = "K139"; Value = 10222.8 },
        @{ Cell = "M139"; Value = -5082.799999999999 }
    )
    "GSM" = @(
        @{ Cell = "H34"; Value = 0 },
        @{ Cell = "J34"; Value = 0 },
        @{ Cell = "L34"; Value = 0 },
        @{ Cell = "N34"; Value = $null },
        @{ Cell = "H43"; Value = 30548.666 },
        @{ Cell = "I43"; Value = 3999.5 },
        @{ Cell = "J43"; Value = 43823.25 },
        @{ Cell = "K43"; Value = 3999.5 },
        @{ Cell = "L43"; Value = 43823.25 },
        @{ Cell = "M43"; Value = -3848.5 },
        @{ Cell = "N43"; Value = -44125.25 },
        @{ Cell = "H46"; Value = 49973.5 },
        @{ Cell = "I46"; Value = 0 },
        @{ Cell = "J46"; Value = 49973.5 },
        @{ Cell = "K46"; Value = 0 },
        @{ Cell = "L46"; Value = 49973.5 },
        @{ Cell = "M46"; Value = $null },
        @{ Cell = "N46"; Value = -50285.5 },
        @{ Cell = "H70"; Value = 35719680 },
        @{ Cell = "I70"; Value = 83337336 },
        @{ Cell = "J70"; Value = 6438.25 },
        @{ Cell = "K70"; Value = 83337336 },
        @{ Cell = "L70"; Value = 6438.25 },
        @{ Cell = "M70"; Value = -83337066 },
        @{ Cell = "N70"; Value = -6978.25 },
        @{ Cell = "H73"; Value = 35719680 },
        @{ Cell = "I73"; Value = 83337336 },
        @{ Cell = "J73"; Value = 6438.25 },
        @{ Cell = "K73"; Value = 83337336 },
        @{ Cell = "L73"; Value = 6438.25 },
        @{ Cell = "M73"; Value = -83336400 },
        @{ Cell = "N73"; Value = -8310.25 },
        @{ Cell = "H76"; Value = 0 },
        @{ Cell = "J76"; Value = 0 },
        @{ Cell = "L76"; Value = 0 },
        @{ Cell = "N76"; Value = $null },
        @{ Cell = "H79"; Value = 0 },
        @{ Cell = "J79"; Value = 0 },
        @{ Cell = "L79"; Value = 0 },
        @{ Cell = "N79"; Value = $null },
        @{ Cell = "H122"; Value = 5132719.5 },
        @{ Cell = "I122"; Value = 7696076.5 },
        @{ Cell = "J122"; Value = 6006.2 },
        @{ Cell = "K122"; Value = 23088229.5 },
        @{ Cell = "L122"; Value = 18018.6 },
        @{ Cell = "M122"; Value = -23085779.5 },
        @{ Cell = "N122"; Value = -22918.6 },
        @{ Cell = "H126"; Value = 7150.706 },
        @{ Cell = "I126"; Value = 1764.5555 },
        @{ Cell = "J126"; Value = 13210.125 },
        @{ Cell = "K126"; Value = 5293.666499999999 },
        @{ Cell = "L126"; Value = 39630.375 },
        @{ Cell = "M126"; Value = -2823.666499999999 },
        @{ Cell = "N126"; Value = -44570.375 },
        @{ Cell = "H132"; Value = 3202.4 },
        @{ Cell = "I132"; Value = 3004 },
        @{ Cell = "J132"; Value = 3500 },
        @{ Cell = "K132"; Value = 9012 },
        @{ Cell = "M132"; Value = -6482 },
        @{ Cell = "N132"; Value = -15560 },
        @{ Cell = "H133"; Value = 96439.125 },
        @{ Cell = "I133"; Value = 0 },
        @{ Cell = "J133"; Value = 96439.125 },
        @{ Cell = "K133"; Value = 0 },
        @{ Cell = "L133"; Value = 96439.125 },
        @{ Cell = "M133"; Value = $null },
        @{ Cell = "N133"; Value = -106559.125 },
        @{ Cell = "H141"; Value = 0 },
        @{ Cell = "J141"; Value = 0 },
        @{ Cell = "L141"; Value = 0 },
        @{ Cell = "N141"; Value = $null }
    )
    "LTW" = @(
        @{ Cell = "H2"; Value = 33394666 },
        @{ Cell = "J2"; Value = 33394666 },
        @{ Cell = "L2"; Value = 33394666 },
        @{ Cell = "N2"; Value = -33394890 },
        @{ Cell = "H40"; Value = 13498.167 },
        @{ Cell = "I40"; Value = 15597.8 },
        @{ Cell = "K40"; Value = 15597.8 },
        @{ Cell = "M40"; Value = -15461.8 },
        @{ Cell = "H46"; Value = 4559.1816 },
        @{ Cell = "I46"; Value = 3518.875 },
        @{ Cell = "K46"; Value = 3518.875 },
        @{ Cell = "M46"; Value = -3330.875 },
        @{ Cell = "H55"; Value = 457.59375 },
        @{ Cell = "I55"; Value = 326.3 },
        @{ Cell = "J55"; Value = 676.4167 },
        @{ Cell = "K55"; Value = 326.3 },
        @{ Cell = "L55"; Value = 676.4167 },
        @{ Cell = "M55"; Value = -153.3 },
        @{ Cell = "N55"; Value = -1022.4167 },
        @{ Cell = "H104"; Value = 39665 },
        @{ Cell = "J104"; Value = 39665 },
        @{ Cell = "L104"; Value = 39665 },
        @{ Cell = "N104"; Value = -46653 },
        @{ Cell = "H122"; Value = 14074.889 },
        @{ Cell = "I122"; Value = 11333 },
        @{ Cell = "J122"; Value = 15445.833 },
        @{ Cell = "K122"; Value = 33999 },
        @{ Cell = "L122"; Value = 46337.499 },
        @{ Cell = "M122"; Value = -31549 },
        @{ Cell = "N122"; Value = -51237.499 },
        @{ Cell = "H132"; Value = 6205 },
        @{ Cell = "I132"; Value = 4071.85 },
        @{ Cell = "J132"; Value = 13315.5 },
        @{ Cell = "K132"; Value = 12215.55 },
        @{ Cell = "L132"; Value = 39946.5 },
        @{ Cell = "M132"; Value = -9685.549999999999 },
        @{ Cell = "N132"; Value = -45006.5 },
        @{ Cell = "H136"; Value = 5231.647 },
        @{ Cell = "I136"; Value = 3808.625 },
        @{ Cell = "J136"; Value = 6496.5557 },
        @{ Cell = "K136"; Value = 11425.875 },
        @{ Cell = "L136"; Value = 19489.6671 },
        @{ Cell = "M136"; Value = -8875.875 },
        @{ Cell = "N136"; Value = -24589.6671 }
    )
    "WVR" = @(
        @{ Cell = "H21"; Value = 20007.5 },
        @{ Cell = "I21"; Value = 20007.5 },
        @{ Cell = "K21"; Value = 20007.5 },
        @{ Cell = "M21"; Value = -19772.5 },
        @{ Cell = "H35"; Value = 20007.5 },
        @{ Cell = "I35"; Value = 20007.5 },
        @{ Cell = "K35"; Value = 20007.5 },
        @{ Cell = "M35"; Value = -19717.5 },
        @{ Cell = "H39"; Value = 24044 },
        @{ Cell = "J39"; Value = 0 },
        @{ Cell = "L39"; Value = 0 },
        @{ Cell = "N39"; Value = $null },
        @{ Cell = "H41"; Value = 17758.143 },
        @{ Cell = "J41"; Value = 17307.8 },
        @{ Cell = "L41"; Value = 17307.8 },
        @{ Cell = "N41"; Value = -18087.8 },
        @{ Cell = "H43"; Value = 22513.5 },
        @{ Cell = "I43"; Value = 15027 },
        @{ Cell = "J43"; Value = 30000 },
        @{ Cell = "K43"; Value = 15027 },
        @{ Cell = "L43"; Value = 30000 },
        @{ Cell = "M43"; Value = -14878 },
        @{ Cell = "N43"; Value = -30298 },
        @{ Cell = "H81"; Value = 4400.4 },
        @{ Cell = "I81"; Value = 3714.9285 },
        @{ Cell = "K81"; Value = 7429.857 },
        @{ Cell = "M81"; Value = -6368.857 },
        @{ Cell = "H84"; Value = 4400.4 },
        @{ Cell = "I84"; Value = 3714.9285 },
        @{ Cell = "K84"; Value = 37149.285 },
        @{ Cell = "M84"; Value = -31845.285 },
        @{ Cell = "H96"; Value = 1666.6666 },
        @{ Cell = "J96"; Value = 1100 },
        @{ Cell = "L96"; Value = 1100 },
        @{ Cell = "N96"; Value = -3846 },
        @{ Cell = "H122"; Value = 11367334 },
        @{ Cell = "I122"; Value = 4101.9443 },
        @{ Cell = "K122"; Value = 12305.8329 },
        @{ Cell = "M122"; Value = -9855.832900000001 },
        @{ Cell = "H126"; Value = 3249.6667 },
        @{ Cell = "I126"; Value = 3199.6 },
        @{ Cell = "K126"; Value = 9598.799999999999 },
        @{ Cell = "M126"; Value = -7128.799999999999 },
        @{ Cell = "H132"; Value = 1628.75 },
        @{ Cell = "I132"; Value = 1518.7059 },
        @{ Cell = "K132"; Value = 4556.1177 },
        @{ Cell = "M132"; Value = -2026.1177 },
        @{ Cell = "H136"; Value = 557822.0600000001 },
        @{ Cell = "I136"; Value = 3233.2778 },
        @{ Cell = "K136"; Value = 9699.8334 },
        @{ Cell = "M136"; Value = -7149.8334 }
    )
}

foreach ($sheetName in $sheetChanges.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($chg in $sheetChanges[$sheetName]) {
        $ws.Range($chg.Cell).Value = $chg.Value
    }
}

Write-Host "Applied $($sheetChanges.Values | ForEach-Object { $_.Count } | Measure-Object -Sum).Sum changes"